$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. New column AU: "Status as of July 4, 2025"
# ---------------------------------------------------------------------------
$ws.Range("AU1").Value = "Status as of July 4, 2025"

# ---------------------------------------------------------------------------
# 2. Drop the stray empty cells that used to hold empty inline strings.
# ---------------------------------------------------------------------------
$emptyCells = @(
    "S2","T2","U2","V2","W2","X2","AA2","AN2","AP2","AQ2","AR2","AS2",
    "R3","S3","T3","AN3","AO3","AQ3","AR3",
    "Q4","R4","T4","Y4","AN4","AO4","AQ4",
    "Q5","R5","S5","T5","Y5","Z5","AN5","AO5","AQ5"
)
foreach ($addr in $emptyCells) {
    $ws.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# 3. Add the (hidden) DropdownOptions sheet right after Sheet1 and fill it
#    with the percentage buckets used for the new status dropdown.
# ---------------------------------------------------------------------------
$ddSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ddSheet.Name = "DropdownOptions"

$optRange = $ddSheet.Range("A1:A7")
$optRange.NumberFormat = "@"
$options = @("0% - 10%", "11% - 25%", "26% - 50%", "51% - 75%", "76% - 90%", "91% - 99%", "100%")
for ($i = 0; $i -lt $options.Length; $i++) {
    $ddSheet.Cells.Item($i + 1, 1).Value = $options[$i]
}
$optRange.Style = "Normal"

# ---------------------------------------------------------------------------
# 4. Data validation dropdown on AU2:AU5 pointing at DropdownOptions!A1:A7
# ---------------------------------------------------------------------------
$target = $ws.Range("AU2:AU5")
$target.Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$target.Validation.InCellDropdown = $true
$target.Validation.ShowInput = $false
$target.Validation.ShowError = $false

# ---------------------------------------------------------------------------
# 5. Keep DropdownOptions hidden, matching the original workbook.
# ---------------------------------------------------------------------------
$ddSheet.Visible = $false
